# Update countries & provincias Spain
#
# 1) Refresh the raw case numbers for the countries whose source data
#    changed (Polonia, Rumania, Dinamarca, Australia, Malasia, Oman,
#    Camerun, Eslovenia, Eritrea).
# 2) Re-sort the whole country table (rows 4-216) by "Casos totales"
#    (column B) descending - this naturally reshuffles Polonia above
#    Corea del Sur, Oman above Islandia and Camerun above Bosnia y
#    Herzegovina, matching the new totals.
# 3) Bump the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update raw numbers (Casos totales, Nuevos casos, Casos activos,
#        Recuperados, Casos criticos, Muertes hoy, Muertes) -------------

# Polonia (row 34)
$ws.Range("B34").Value = 10759
$ws.Range("C34").Value = 248
$ws.Range("D34").Value = 1944
$ws.Range("E34").Value = 8352
$ws.Range("F34").Value = 160
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 463

# Rumania (row 35)
$ws.Range("B35").Value = 10096
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 2478
$ws.Range("E35").Value = 7066
$ws.Range("F35").Value = 236
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 552

# Dinamarca (row 37)
$ws.Range("B37").Value = 8210
$ws.Range("C37").Value = 137
$ws.Range("D37").Value = 5384
$ws.Range("E37").Value = 2432
$ws.Range("F37").Value = 74
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 394

# Australia (row 46)
$ws.Range("B46").Value = 6675
$ws.Range("C46").Value = 8
$ws.Range("D46").Value = 5136
$ws.Range("E46").Value = 1460
$ws.Range("F46").Value = 43
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 79

# Malasia (row 47)
$ws.Range("B47").Value = 5691
$ws.Range("C47").Value = 88
$ws.Range("D47").Value = 3663
$ws.Range("E47").Value = 1932
$ws.Range("F47").Value = 41
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 96

# Oman (row 69)
$ws.Range("B69").Value = 1790
$ws.Range("C69").Value = 74
$ws.Range("D69").Value = 325
$ws.Range("E69").Value = 1456
$ws.Range("F69").Value = 3
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 9

# Eslovenia (row 77)
$ws.Range("B77").Value = 1360
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 211
$ws.Range("E77").Value = 1076
$ws.Range("F77").Value = 23
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 79

# Camerun (row 79)
$ws.Range("B79").Value = 1430
$ws.Range("C79").Value = 96
$ws.Range("D79").Value = 668
$ws.Range("E79").Value = 719
$ws.Range("F79").Value = 20
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 43

# Eritrea (row 170)
$ws.Range("D170").Value = 11
$ws.Range("E170").Value = 28

# --- 2) Re-sort the country table by Casos totales (column B), desc ----

$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)

# --- 3) Bump "last updated" timestamp ----------------------------------

$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 10:52"
